# Optimized __get_list_of_games: add a third worksheet "Opt2" run with
# updated listOfGames timings, and tidy the sheet names (drop the leading
# "After " prefix now that all tabs in the book represent "after" runs).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename the two existing sheets ------------------------------------
$ws1.Name = " _check_year opt"
$ws2.Name = "listOfGames Opt1"

# --- Update the D column (third listOfGames timing run) on sheet 2 -----
$ws2.Range("D2").Value = 34.156999999999996
$ws2.Range("D3").Value = 33.08
$ws2.Range("D4").Value = 1.07

# --- Show the AVERAGE column rounded to 2 decimal places ---------------
$ws2.Range("E1").NumberFormat = "0.00"
$ws2.Range("E2:E6").NumberFormat = "0.00"

# --- Duplicate the (now-updated) sheet to capture the Opt2 results -----
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = " listOfGames Opt2"

# --- Restore sensible selections / active sheet -------------------------
[void]$ws2.Range("A1:E6").Select()
$ws3.Activate()
[void]$ws3.Range("A1:E6").Select()
